{"js": "// Replace the 25 \"A\u00d7B=C\" answer cells in the practice-problems table\n// with their updated values, per the commit diff. Each old value is\n// unique within the document, so a scoped search-and-replace on the\n// document body is safe and keeps every run's original formatting\n// (font/size) untouched -- only the text content changes.\nconst replacements = [\n  [\"691\u00d72=1382\", \"525\u00d72=1050\"],\n  [\"764\u00d73=2292\", \"435\u00d77=3045\"],\n  [\"835\u00d73=2505\", \"373\u00d78=2984\"],\n  [\"901\u00d76=5406\", \"745\u00d78=5960\"],\n  [\"300\u00d73=900\", \"490\u00d73=1470\"],\n  [\"778\u00d76=4668\", \"503\u00d77=3521\"],\n  [\"897\u00d73=2691\", \"875\u00d79=7875\"],\n  [\"322\u00d79=2898\", \"593\u00d74=2372\"],\n  [\"383\u00d72=766\", \"582\u00d79=5238\"],\n  [\"515\u00d72=1030\", \"227\u00d74=908\"],\n  [\"121\u00d73=363\", \"946\u00d76=5676\"],\n  [\"975\u00d75=4875\", \"101\u00d77=707\"],\n  [\"151\u00d76=906\", \"425\u00d77=2975\"],\n  [\"251\u00d73=753\", \"451\u00d79=4059\"],\n  [\"251\u00d72=502\", \"482\u00d79=4338\"],\n  [\"831\u00d72=1662\", \"128\u00d75=640\"],\n  [\"370\u00d79=3330\", \"776\u00d75=3880\"],\n  [\"802\u00d74=3208\", \"889\u00d74=3556\"],\n  [\"421\u00d73=1263\", \"917\u00d77=6419\"],\n  [\"822\u00d74=3288\", \"665\u00d79=5985\"],\n  [\"177\u00d72=354\", \"807\u00d75=4035\"],\n  [\"868\u00d77=6076\", \"987\u00d75=4935\"],\n  [\"183\u00d77=1281\", \"736\u00d76=4416\"],\n  [\"763\u00d76=4578\", \"247\u00d74=988\"],\n  [\"910\u00d79=8190\", \"689\u00d75=3445\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"A\u00d7B=C\" answer cells in the practice-problems table\n# with their updated values, per the commit diff. Each old value is\n# unique within the document, so a Find/Replace over the whole document\n# Content range is safe and leaves every run's original formatting\n# (font/size) untouched -- only the text content changes.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"691\u00d72=1382\"; New = \"525\u00d72=1050\" },\n    @{ Old = \"764\u00d73=2292\"; New = \"435\u00d77=3045\" },\n    @{ Old = \"835\u00d73=2505\"; New = \"373\u00d78=2984\" },\n    @{ Old = \"901\u00d76=5406\"; New = \"745\u00d78=5960\" },\n    @{ Old = \"300\u00d73=900\"; New = \"490\u00d73=1470\" },\n    @{ Old = \"778\u00d76=4668\"; New = \"503\u00d77=3521\" },\n    @{ Old = \"897\u00d73=2691\"; New = \"875\u00d79=7875\" },\n    @{ Old = \"322\u00d79=2898\"; New = \"593\u00d74=2372\" },\n    @{ Old = \"383\u00d72=766\"; New = \"582\u00d79=5238\" },\n    @{ Old = \"515\u00d72=1030\"; New = \"227\u00d74=908\" },\n    @{ Old = \"121\u00d73=363\"; New = \"946\u00d76=5676\" },\n    @{ Old = \"975\u00d75=4875\"; New = \"101\u00d77=707\" },\n    @{ Old = \"151\u00d76=906\"; New = \"425\u00d77=2975\" },\n    @{ Old = \"251\u00d73=753\"; New = \"451\u00d79=4059\" },\n    @{ Old = \"251\u00d72=502\"; New = \"482\u00d79=4338\" },\n    @{ Old = \"831\u00d72=1662\"; New = \"128\u00d75=640\" },\n    @{ Old = \"370\u00d79=3330\"; New = \"776\u00d75=3880\" },\n    @{ Old = \"802\u00d74=3208\"; New = \"889\u00d74=3556\" },\n    @{ Old = \"421\u00d73=1263\"; New = \"917\u00d77=6419\" },\n    @{ Old = \"822\u00d74=3288\"; New = \"665\u00d79=5985\" },\n    @{ Old = \"177\u00d72=354\"; New = \"807\u00d75=4035\" },\n    @{ Old = \"868\u00d77=6076\"; New = \"987\u00d75=4935\" },\n    @{ Old = \"183\u00d77=1281\"; New = \"736\u00d76=4416\" },\n    @{ Old = \"763\u00d76=4578\"; New = \"247\u00d74=988\" },\n    @{ Old = \"910\u00d79=8190\"; New = \"689\u00d75=3445\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $result = $range.Find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n\n    if (-not $result) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
